# Adds the new student "Xulia Teixeira " as row 29 on every sheet, fills in
# the previously-empty attitude-score columns (B on sheet1-3, C on sheet4)
# for every existing student, and updates the active sheet/selection to
# match the editor's final position (Ninja sheet, cell C29).

$wb = $excel.ActiveWorkbook

$newName = "Xulia Teixeira "

# ---------------------------------------------------------------------
# Sheet "Astronauta" (sheet1): column B, rows 2-28 + new row 29
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Astronauta")

$bVals = @(0,0,0,1,0,0,0,1,0,0,0,0,0,0,0,1,1,0,1,0,0,1,1,1,1,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Cells.Item($r, 2)
    if ($r -eq 3) {
        # B3 previously carried an (unused) underline style - clear it
        # along with the value so the cell becomes a plain number.
        $cell.ClearFormats()
    }
    $cell.Value = $bVals[$i]
}

$ws.Cells.Item(29, 1).Value = $newName
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 12).FormulaR1C1 = "=IFERROR(SUM(RC[-10]:RC[-1])/COUNT(RC[-10]:RC[-1])*100, 0)"

# ---------------------------------------------------------------------
# Sheet "Senador" (sheet2): column B, rows 2-28 + new row 29 (no L29)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Senador")

$bVals = @(0,0,1,1,0,0,0,0,1,1,0,0,0,0,0,1,1,0,1,1,0,1,1,1,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
}

$ws.Cells.Item(29, 1).Value = $newName
$ws.Cells.Item(29, 2).Value = 1

# ---------------------------------------------------------------------
# Sheet "Mago" (sheet3): column B, rows 2-28 + new row 29
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mago")

$bVals = @(0,0,1,1,0,0,1,1,1,1,0,0,1,0,1,1,1,1,1,1,1,1,1,1,1,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
}

$ws.Cells.Item(29, 1).Value = $newName
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 12).FormulaR1C1 = "=IFERROR(SUM(RC[-10]:RC[-1])/COUNT(RC[-10]:RC[-1])*100, 0)"

# ---------------------------------------------------------------------
# Sheet "Ninja" (sheet4): column C, rows 2-28 + new row 29 (B & C)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ninja")

$cVals = @(0,1,1,1,1,1,1,1,1,1,1,0,1,1,0,1,1,1,1,1,0,1,1,1,0,0,0)
for ($i = 0; $i -lt $cVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
}

$ws.Cells.Item(29, 1).Value = $newName
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 16).FormulaR1C1 = "=IFERROR(SUM(RC[-14]:RC[-1])/COUNT(RC[-14]:RC[-1])*100, 0)"

# ---------------------------------------------------------------------
# Final view state: Ninja is the active sheet, with the selections the
# editor left on each sheet.
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Astronauta")
$wsA.Activate() | Out-Null
$wsA.Range("C31").Select() | Out-Null

$wsS = $wb.Worksheets.Item("Senador")
$wsS.Activate() | Out-Null
$wsS.Range("B30").Select() | Out-Null

$wsM = $wb.Worksheets.Item("Mago")
$wsM.Activate() | Out-Null
$wsM.Range("B29").Select() | Out-Null

$wsN = $wb.Worksheets.Item("Ninja")
$wsN.Activate() | Out-Null
$wsN.Range("C29").Select() | Out-Null
